# Horarios actualizados Línea 141 - 1220
# Scrape refresh: "Última actualización" 04:43:39 -> 04:56:49, with updated
# "Minutos" countdowns and newly-arrived rows appended at the bottom of the
# sheets that gained stops.

$wb = $excel.ActiveWorkbook

$oldStamp = "04:43:39"
$newStamp = "04:56:49"

# ---------------------------------------------------------------------
# Sheet 1: LP1912  (14 filas -> 16 filas; 2 new rows appended at 20/21)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: " + $newStamp
$ws1.Range("A3").Value = "Total filas: 16"

$sheet1Minutes = @{
    8  = 20
    9  = 26
    10 = 38
    11 = 50
    12 = 58
    13 = 68
    14 = 75
    15 = 78
    16 = 85
    17 = 91
    18 = 93
    19 = 95
}

foreach ($row in 8..19) {
    $ws1.Cells.Item($row, 1).Value = $newStamp
    $ws1.Cells.Item($row, 4).Value = $sheet1Minutes[$row]
}

# Newly-scraped arrivals appended to the bottom of the sheet.
$ws1.Cells.Item(20, 1).Value = $newStamp
$ws1.Cells.Item(20, 2).Value = "06:44"
$ws1.Cells.Item(20, 3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(20, 4).Value = 108
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = $newStamp
$ws1.Cells.Item(21, 2).Value = "06:46"
$ws1.Cells.Item(21, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(21, 4).Value = 110
$ws1.Cells.Item(21, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215  (3 filas -> 4 filas; 1 new row appended at 9)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: " + $newStamp
$ws2.Range("A3").Value = "Total filas: 4"

$sheet2Minutes = @{
    7 = 38
    8 = 75
}

foreach ($row in 7..8) {
    $ws2.Cells.Item($row, 1).Value = $newStamp
    $ws2.Cells.Item($row, 4).Value = $sheet2Minutes[$row]
}

$ws2.Cells.Item(9, 1).Value = $newStamp
$ws2.Cells.Item(9, 2).Value = "06:46"
$ws2.Cells.Item(9, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(9, 4).Value = 110
$ws2.Cells.Item(9, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173  (still 3 filas; only timestamps/minutes refresh)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: " + $newStamp

$sheet3Minutes = @{
    6 = 47
    7 = 72
    8 = 96
}

foreach ($row in 6..8) {
    $ws3.Cells.Item($row, 1).Value = $newStamp
    $ws3.Cells.Item($row, 4).Value = $sheet3Minutes[$row]
}
